$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = 1769400
$ws.Range("E8").Value = 1770800
$ws.Range("F8").Value = 1834400
$ws.Range("G8").Value = 1733600
$ws.Range("H8").Value = 1751800
$ws.Range("I8").Value = 1629300
$ws.Range("J8").Value = 1553900

$ws.Range("D9").Value = 830900
$ws.Range("E9").Value = 840300
$ws.Range("F9").Value = 866900
$ws.Range("G9").Value = 820100
$ws.Range("H9").Value = 822700
$ws.Range("I9").Value = 764300
$ws.Range("J9").Value = 740300

$ws.Range("D10").Value = 938500
$ws.Range("E10").Value = 930500
$ws.Range("F10").Value = 967400
$ws.Range("G10").Value = 913400
$ws.Range("H10").Value = 929100
$ws.Range("I10").Value = 865000
$ws.Range("J10").Value = 813700

$ws.Range("D12").Value = 7100
$ws.Range("E12").Value = 7300
$ws.Range("F12").Value = 7600
$ws.Range("G12").Value = 7900
$ws.Range("H12").Value = 7300
$ws.Range("I12").Value = 7100
$ws.Range("J12").Value = 7200

$ws.Range("D14").Value = 1900
$ws.Range("G14").Value = 54400
$ws.Range("I14").Value = 26700
$ws.Range("J14").Value = 7500

$ws.Range("D17").Value = 1656000
$ws.Range("E17").Value = 1670700
$ws.Range("F17").Value = 1709000
$ws.Range("G17").Value = 1669500
$ws.Range("H17").Value = 1626500
$ws.Range("I17").Value = 1552400
$ws.Range("J17").Value = 1467700

$ws.Range("D18").Value = 113300
$ws.Range("E18").Value = 100000
$ws.Range("F18").Value = 125300
$ws.Range("G18").Value = 64000
$ws.Range("H18").Value = 125300
$ws.Range("I18").Value = 76800
$ws.Range("J18").Value = 86300

$ws.Range("D20").Value = 15900
$ws.Range("E20").Value = 50000
$ws.Range("F20").Value = 10500
$ws.Range("G20").Value = 39400
$ws.Range("H20").Value = 11600
$ws.Range("I20").Value = 23500
$ws.Range("J20").Value = 6800

$ws.Range("D21").Value = 178900
$ws.Range("E21").Value = 195600
$ws.Range("F21").Value = 179400
$ws.Range("G21").Value = 149300
$ws.Range("H21").Value = 182400
$ws.Range("I21").Value = 144600
$ws.Range("J21").Value = 135300

$ws.Range("D23").Value = 129100
$ws.Range("E23").Value = 149800
$ws.Range("F23").Value = 135200
$ws.Range("G23").Value = 102500
$ws.Range("H23").Value = 135900
$ws.Range("I23").Value = 98900
$ws.Range("J23").Value = 92300

$ws.Range("D24").Value = 47500
$ws.Range("E24").Value = 47700
$ws.Range("F24").Value = 42800
$ws.Range("G24").Value = 29900
$ws.Range("H24").Value = 51000
$ws.Range("I24").Value = 34100
$ws.Range("J24").Value = 38000

$ws.Range("D26").Value = 81700
$ws.Range("E26").Value = 102100
$ws.Range("F26").Value = 92500
$ws.Range("G26").Value = 72700
$ws.Range("H26").Value = 84900
$ws.Range("I26").Value = 64800
$ws.Range("J26").Value = 54300

$ws.Range("D27").Value = 90700
$ws.Range("E27").Value = 113200
$ws.Range("F27").Value = 100900
$ws.Range("G27").Value = 76300
$ws.Range("H27").Value = 91400
$ws.Range("I27").Value = 71200
$ws.Range("J27").Value = 62500

$ws.Range("D32").Value = -15900
$ws.Range("E32").Value = -50000
$ws.Range("F32").Value = -10500
$ws.Range("G32").Value = -39400
$ws.Range("H32").Value = -11600
$ws.Range("I32").Value = -23500
$ws.Range("J32").Value = -6800

$ws.Range("D33").Value = 88100
$ws.Range("E33").Value = 113200
$ws.Range("F33").Value = 100900
$ws.Range("G33").Value = 76300
$ws.Range("H33").Value = 91400
$ws.Range("I33").Value = 71200
$ws.Range("J33").Value = 62500

$ws.Range("D35").Value = 88100
$ws.Range("E35").Value = 113200
$ws.Range("F35").Value = 100900
$ws.Range("G35").Value = 76300
$ws.Range("H35").Value = 91400
$ws.Range("I35").Value = 71200
$ws.Range("J35").Value = 62500

$ws.Range("D41").Value = 305400
$ws.Range("E41").Value = 331900
$ws.Range("F41").Value = 327200
$ws.Range("G41").Value = 371500
$ws.Range("H41").Value = 296700
$ws.Range("I41").Value = 238900
$ws.Range("J41").Value = 277700

$ws.Range("D42").Value = 14200
$ws.Range("E42").Value = 13200
$ws.Range("F42").Value = 17000
$ws.Range("G42").Value = 21600
$ws.Range("H42").Value = 31800
$ws.Range("I42").Value = 41600
$ws.Range("J42").Value = 46800

$ws.Range("D43").Value = 211700
$ws.Range("E43").Value = 208700
$ws.Range("F43").Value = 223400
$ws.Range("G43").Value = 218200
$ws.Range("H43").Value = 216500
$ws.Range("I43").Value = 199400
$ws.Range("J43").Value = 192200

$ws.Range("D44").Value = 385800
$ws.Range("E44").Value = 396200
$ws.Range("F44").Value = 401800
$ws.Range("G44").Value = 387800
$ws.Range("H44").Value = 363500
$ws.Range("I44").Value = 342600
$ws.Range("J44").Value = 296900

$ws.Range("D45").Value = 39500
$ws.Range("E45").Value = 78900
$ws.Range("F45").Value = 87000
$ws.Range("G45").Value = 88200
$ws.Range("H45").Value = 81200
$ws.Range("I45").Value = 113900
$ws.Range("J45").Value = 65900

$ws.Range("D46").Value = 956500
$ws.Range("E46").Value = 1028900
$ws.Range("F46").Value = 1056300
$ws.Range("G46").Value = 1087200
$ws.Range("H46").Value = 989800
$ws.Range("I46").Value = 936400
$ws.Range("J46").Value = 879500

$ws.Range("D47").Value = 830100
$ws.Range("E47").Value = 729700
$ws.Range("F47").Value = 693700
$ws.Range("G47").Value = 741400
$ws.Range("H47").Value = 586200
$ws.Range("I47").Value = 542100
$ws.Range("J47").Value = 439900

$ws.Range("D48").Value = 491200
$ws.Range("E48").Value = 499800
$ws.Range("F48").Value = 487600
$ws.Range("G48").Value = 444700
$ws.Range("H48").Value = 442800
$ws.Range("I48").Value = 449000
$ws.Range("J48").Value = 443700

$ws.Range("D49").Value = 266300
$ws.Range("E49").Value = 252400
$ws.Range("F49").Value = 271400
$ws.Range("G49").Value = 284700
$ws.Range("H49").Value = 329200
$ws.Range("I49").Value = 301400
$ws.Range("J49").Value = 180000

$ws.Range("D52").Value = 154600
$ws.Range("E52").Value = 155700
$ws.Range("F52").Value = 138400
$ws.Range("G52").Value = 156500
$ws.Range("H52").Value = 110900
$ws.Range("I52").Value = 72100
$ws.Range("J52").Value = 55600

$ws.Range("D54").Value = 2698700
$ws.Range("E54").Value = 2666400
$ws.Range("F54").Value = 2647400
$ws.Range("G54").Value = 2714500
$ws.Range("H54").Value = 2458800
$ws.Range("I54").Value = 2301000
$ws.Range("J54").Value = 1998700

$ws.Range("D57").Value = 113600
$ws.Range("E57").Value = 117900
$ws.Range("F57").Value = 121600
$ws.Range("G57").Value = 111900
$ws.Range("H57").Value = 106000
$ws.Range("I57").Value = 112100
$ws.Range("J57").Value = 97100

$ws.Range("D58").Value = 64700
$ws.Range("E58").Value = 70200
$ws.Range("F58").Value = 109000
$ws.Range("G58").Value = 95600
$ws.Range("H58").Value = 156600
$ws.Range("I58").Value = 155900
$ws.Range("J58").Value = 65700

$ws.Range("D59").Value = 180000
$ws.Range("E59").Value = 183100
$ws.Range("F59").Value = 170300
$ws.Range("G59").Value = 175000
$ws.Range("H59").Value = 158800
$ws.Range("I59").Value = 188500
$ws.Range("J59").Value = 159100

$ws.Range("D60").Value = 358300
$ws.Range("E60").Value = 371200
$ws.Range("F60").Value = 400900
$ws.Range("G60").Value = 382400
$ws.Range("H60").Value = 421300
$ws.Range("I60").Value = 456600
$ws.Range("J60").Value = 321900

$ws.Range("E61").Value = 1700
$ws.Range("F61").Value = 900
$ws.Range("G61").Value = 38400
$ws.Range("H61").Value = 7400
$ws.Range("I61").Value = 13700
$ws.Range("J61").Value = 5800

$ws.Range("D62").Value = 192300
$ws.Range("E62").Value = 191900
$ws.Range("F62").Value = 171900
$ws.Range("G62").Value = 200300
$ws.Range("H62").Value = 153900
$ws.Range("I62").Value = 123800
$ws.Range("J62").Value = 103300

$ws.Range("D66").Value = 595000
$ws.Range("E66").Value = 609200
$ws.Range("F66").Value = 619100
$ws.Range("G66").Value = 645600
$ws.Range("H66").Value = 604600
$ws.Range("I66").Value = 613700
$ws.Range("J66").Value = 448400

$ws.Range("D72").Value = 1558700
$ws.Range("E72").Value = 1537400
$ws.Range("F72").Value = 1466300
$ws.Range("G72").Value = 1403600
$ws.Range("H72").Value = 1369300
$ws.Range("I72").Value = 1313600
$ws.Range("J72").Value = 1278000

$ws.Range("D76").Value = 2103700
$ws.Range("E76").Value = 2057200
$ws.Range("F76").Value = 2028300
$ws.Range("G76").Value = 2068900
$ws.Range("H76").Value = 1854200
$ws.Range("I76").Value = 1687300
$ws.Range("J76").Value = 1550300

$ws.Range("D81").Value = 88100
$ws.Range("E81").Value = 113200
$ws.Range("F81").Value = 100900
$ws.Range("G81").Value = 76300
$ws.Range("H81").Value = 91400
$ws.Range("I81").Value = 71200
$ws.Range("J81").Value = 62500

$ws.Range("D83").Value = 49600
$ws.Range("E83").Value = 45500
$ws.Range("F83").Value = 43500
$ws.Range("G83").Value = 45900
$ws.Range("H83").Value = 45500
$ws.Range("I83").Value = 44200
$ws.Range("J83").Value = 42100

$ws.Range("D89").Value = 140100
$ws.Range("E89").Value = 147800
$ws.Range("F89").Value = 114200
$ws.Range("G89").Value = 129600
$ws.Range("H89").Value = 80900
$ws.Range("I89").Value = 111300
$ws.Range("J89").Value = 90900

$ws.Range("D91").Value = -31000
$ws.Range("E91").Value = -49800
$ws.Range("F91").Value = -68200
$ws.Range("G91").Value = -28000
$ws.Range("H91").Value = -20500
$ws.Range("I91").Value = -22400
$ws.Range("J91").Value = -24500

$ws.Range("D94").Value = -66600
$ws.Range("E94").Value = -27400
$ws.Range("F94").Value = -103100
$ws.Range("H94").Value = 15000
$ws.Range("I94").Value = -212600
$ws.Range("J94").Value = -31300

$ws.Range("D96").Value = -66800
$ws.Range("E96").Value = -42000
$ws.Range("F96").Value = -38200
$ws.Range("G96").Value = -42000
$ws.Range("H96").Value = -35700
$ws.Range("I96").Value = -35700
$ws.Range("J96").Value = -25500

$ws.Range("D100").Value = -111200
$ws.Range("E100").Value = -118000
$ws.Range("F100").Value = -41100
$ws.Range("G100").Value = -75900
$ws.Range("H100").Value = -50200
$ws.Range("I100").Value = 48600
$ws.Range("J100").Value = -25500

$ws.Range("E101").Value = -3000
$ws.Range("F101").Value = -9300
$ws.Range("G101").Value = 14800
$ws.Range("H101").Value = 9900
$ws.Range("I101").Value = 3300

$ws.Range("D102").Value = -40800
$ws.Range("F102").Value = -39300
$ws.Range("G102").Value = 70100
$ws.Range("H102").Value = 55500
$ws.Range("I102").Value = -49500
$ws.Range("J102").Value = 33200
